# Blueprint "Individuals in need" update
# - refresh wording around the DONOTHROW association branding
# - refresh notification / food-choice / departure-reason copy
# - clear the stray "new reason of departure" DB note in row 14
# - shrink rows 9 & 11 (text got shorter) and move the viewport/selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- text content updates -------------------------------------------------

# Row 7 / Web presentation blurb now references "Donothrow" instead of
# "the association"
$ws.Range("E7").Value  = "Presentation of Donothrow"

# Row 9 (Web channel) & row 11 (Mobile channel) share the same wording
$ws.Range("D9").Value  = "Homepage of the website with information on Donothrow"
$ws.Range("F9").Value  = "Presentation of Donothrow, link to the website"
$ws.Range("F11").Value = "Presentation of Donothrow, link to the website"

$ws.Range("J9").Value  = "Inserts user info and the type of food needed"
$ws.Range("J11").Value = "Inserts user info and the type of food needed"

$ws.Range("K9").Value  = "Reads the notification (if the request was accepted, what food is available)"
$ws.Range("K11").Value = "Reads the notification (if the request was accepted, what food is available)"

$ws.Range("L9").Value  = "Chooses some food from the proposed list"
$ws.Range("L11").Value = "Chooses some food from the proposed list"

$ws.Range("M9").Value  = "Checks some reasons of departure in a list"
$ws.Range("M11").Value = "Checks some reasons of departure in a list"

# Row 14 backstage note about a new "reason of departure" DB row is dropped
$ws.Range("M14").Value = ""
$ws.Range("M14").ClearContents()

# --- row sizing (shorter copy needs less height) ---------------------------
$ws.Rows.Item(9).RowHeight  = 104.4
$ws.Rows.Item(11).RowHeight = 96

# --- viewport / selection ---------------------------------------------------
$ws.Activate()
$ws.Range("E7").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 4
